$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1281.5
$ws.Range("I28").Value = 930.375
$ws.Range("K28").Value = 930.375
$ws.Range("M28").Value = -445.375
$ws.Range("H40").Value = 2850
$ws.Range("H58").Value = 415.22223
$ws.Range("I58").Value = 342.125
$ws.Range("J58").Value = 1000
$ws.Range("K58").Value = 1026.375
$ws.Range("L58").Value = 3000
$ws.Range("M58").Value = -876.375
$ws.Range("N58").Value = -3300
$ws.Range("H111").Value = 4045.0435
$ws.Range("I111").Value = 3914.2942
$ws.Range("J111").Value = 4415.5
$ws.Range("K111").Value = 11742.8826
$ws.Range("L111").Value = 13246.5
$ws.Range("M111").Value = -8675.882599999999
$ws.Range("N111").Value = -19380.5
$ws.Range("H116").Value = 3375.8462
$ws.Range("I116").Value = 3088.6
$ws.Range("K116").Value = 3088.6
$ws.Range("M116").Value = 353.4000000000001
$ws.Range("H118").Value = 1323.75
$ws.Range("I118").Value = 1313.1428
$ws.Range("J118").Value = 1398
$ws.Range("K118").Value = 3939.4284
$ws.Range("L118").Value = 4194
$ws.Range("M118").Value = -2282.4284
$ws.Range("N118").Value = -7508
$ws.Range("H131").Value = 5368.647
$ws.Range("I131").Value = 4140.778
$ws.Range("K131").Value = 12422.334
$ws.Range("M131").Value = -7382.334000000001
$ws.Range("H132").Value = 2707.6765
$ws.Range("I132").Value = 2303.1072
$ws.Range("J132").Value = 4595.6665
$ws.Range("K132").Value = 6909.321599999999
$ws.Range("L132").Value = 13786.9995
$ws.Range("M132").Value = -4379.321599999999
$ws.Range("N132").Value = -18846.9995
$ws.Range("H141").Value = 1672.7273
$ws.Range("I141").Value = 1645
$ws.Range("J141").Value = 1950
$ws.Range("K141").Value = 4935
$ws.Range("L141").Value = 5850
$ws.Range("M141").Value = 245
$ws.Range("N141").Value = -16210

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 758.6111
$ws.Range("I2").Value = 608.9375
$ws.Range("K2").Value = 608.9375
$ws.Range("M2").Value = -495.9375
$ws.Range("H32").Value = 12053162
$ws.Range("I32").Value = 14288304
$ws.Range("J32").Value = 17780.691
$ws.Range("K32").Value = 14288304
$ws.Range("L32").Value = 17780.691
$ws.Range("M32").Value = -14288017
$ws.Range("N32").Value = -18354.691
$ws.Range("H61").Value = 20411328
$ws.Range("I61").Value = 25642810
$ws.Range("K61").Value = 25642810
$ws.Range("M61").Value = -25642598
$ws.Range("H74").Value = 37080172
$ws.Range("I74").Value = 43528236
$ws.Range("K74").Value = 43528236
$ws.Range("M74").Value = -43527362
$ws.Range("H77").Value = 37080172
$ws.Range("I77").Value = 43528236
$ws.Range("K77").Value = 217641180
$ws.Range("M77").Value = -217636812
$ws.Range("H116").Value = 758.6111
$ws.Range("I116").Value = 608.9375
$ws.Range("K116").Value = 608.9375
$ws.Range("M116").Value = 1685.0625
$ws.Range("H132").Value = 20835924
$ws.Range("I132").Value = 2640.3696
$ws.Range("K132").Value = 7921.1088
$ws.Range("M132").Value = -5391.1088
$ws.Range("H136").Value = 20411328
$ws.Range("I136").Value = 25642810
$ws.Range("K136").Value = 76928430
$ws.Range("M136").Value = -76925880

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 758.6111
$ws.Range("I3").Value = 608.9375
$ws.Range("K3").Value = 608.9375
$ws.Range("M3").Value = -494.9375
$ws.Range("H60").Value = 105260.5
$ws.Range("J60").Value = 105260.5
$ws.Range("L60").Value = 105260.5
$ws.Range("N60").Value = -106458.5
$ws.Range("H126").Value = 85000
$ws.Range("J126").Value = 85000
$ws.Range("L126").Value = 85000
$ws.Range("N126").Value = -94880
$ws.Range("H134").Value = 2376.0408
$ws.Range("I134").Value = 2126.1064
$ws.Range("J134").Value = 8249.5
$ws.Range("K134").Value = 6378.3192
$ws.Range("L134").Value = 24748.5
$ws.Range("M134").Value = -3843.3192
$ws.Range("N134").Value = -29818.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 544.75
$ws.Range("I16").Value = 544.75
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 544.75
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -257.75
$ws.Range("H31").Value = 16396720
$ws.Range("I31").Value = 2064.698
$ws.Range("J31").Value = 125011310
$ws.Range("K31").Value = 2064.698
$ws.Range("L31").Value = 125011310
$ws.Range("M31").Value = -1769.698
$ws.Range("N31").Value = -125011900
$ws.Range("H34").Value = 16396720
$ws.Range("I34").Value = 2064.698
$ws.Range("J34").Value = 125011310
$ws.Range("K34").Value = 2064.698
$ws.Range("L34").Value = 125011310
$ws.Range("M34").Value = -1862.698
$ws.Range("N34").Value = -125011714
$ws.Range("H58").Value = 2361.3333
$ws.Range("I58").Value = 1605
$ws.Range("J58").Value = 4198.143
$ws.Range("K58").Value = 1605
$ws.Range("L58").Value = 4198.143
$ws.Range("M58").Value = -1402
$ws.Range("N58").Value = -4604.143
$ws.Range("H105").Value = 8091.579
$ws.Range("I105").Value = 1774.4166
$ws.Range("K105").Value = 1774.4166
$ws.Range("M105").Value = -27.41660000000002
$ws.Range("H107").Value = 2725
$ws.Range("I107").Value = 2000
$ws.Range("J107").Value = 3450
$ws.Range("K107").Value = 2000
$ws.Range("L107").Value = 3450
$ws.Range("M107").Value = -80
$ws.Range("N107").Value = -7290
$ws.Range("H109").Value = 48124.168
$ws.Range("J109").Value = 46749.1
$ws.Range("L109").Value = 46749.1
$ws.Range("N109").Value = -48829.1
$ws.Range("H113").Value = 544.75
$ws.Range("I113").Value = 544.75
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 544.75
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = 1625.25
$ws.Range("H136").Value = 2361.3333
$ws.Range("I136").Value = 1605
$ws.Range("J136").Value = 4198.143
$ws.Range("K136").Value = 4815
$ws.Range("L136").Value = 12594.429
$ws.Range("M136").Value = -2265
$ws.Range("N136").Value = -17694.429

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 5480.85
$ws.Range("I137").Value = 3607.6667
$ws.Range("J137").Value = 6283.643
$ws.Range("K137").Value = 10823.0001
$ws.Range("L137").Value = 18850.929
$ws.Range("M137").Value = -5723.000100000001
$ws.Range("N137").Value = -29050.929

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2917.5186
$ws.Range("I113").Value = 1914.8667
$ws.Range("J113").Value = 4170.8335
$ws.Range("K113").Value = 1914.8667
$ws.Range("L113").Value = 4170.8335
$ws.Range("M113").Value = 255.1333
$ws.Range("N113").Value = -8510.833500000001
$ws.Range("H122").Value = 6298.364
$ws.Range("I122").Value = 3461.5
$ws.Range("J122").Value = 9702.6
$ws.Range("K122").Value = 10384.5
$ws.Range("L122").Value = 29107.8
$ws.Range("M122").Value = -7934.5
$ws.Range("N122").Value = -34007.8
$ws.Range("H132").Value = 2545.3333
$ws.Range("I132").Value = 2152.6072
$ws.Range("J132").Value = 3545
$ws.Range("K132").Value = 6457.821599999999
$ws.Range("L132").Value = 10635
$ws.Range("M132").Value = -3927.821599999999
$ws.Range("N132").Value = -15695
$ws.Range("H137").Value = 104399.2
$ws.Range("I137").Value = 27000
$ws.Range("J137").Value = 123749
$ws.Range("K137").Value = 27000
$ws.Range("L137").Value = 123749
$ws.Range("M137").Value = -21900
$ws.Range("N137").Value = -133949

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 760.55554
$ws.Range("I16").Value = 705.625
$ws.Range("K16").Value = 705.625
$ws.Range("M16").Value = -535.625
$ws.Range("H46").Value = 1262.7333
$ws.Range("I46").Value = 628.36365
$ws.Range("K46").Value = 628.36365
$ws.Range("M46").Value = -440.36365
$ws.Range("H61").Value = 4731.357
$ws.Range("I61").Value = 3917.889
$ws.Range("K61").Value = 3917.889
$ws.Range("M61").Value = -3715.889
$ws.Range("H68").Value = 3008
$ws.Range("I68").Value = 2790
$ws.Range("J68").Value = 3516.6667
$ws.Range("K68").Value = 2790
$ws.Range("L68").Value = 3516.6667
$ws.Range("M68").Value = -2041
$ws.Range("N68").Value = -5014.6667
$ws.Range("H71").Value = 3008
$ws.Range("I71").Value = 2790
$ws.Range("J71").Value = 3516.6667
$ws.Range("K71").Value = 13950
$ws.Range("L71").Value = 17583.3335
$ws.Range("M71").Value = -10206
$ws.Range("N71").Value = -25071.3335
$ws.Range("H113").Value = 4731.357
$ws.Range("I113").Value = 3917.889
$ws.Range("K113").Value = 3917.889
$ws.Range("M113").Value = -1747.889

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 626334.2
$ws.Range("I100").Value = 715524.8
$ws.Range("K100").Value = 1431049.6
$ws.Range("M100").Value = -1430508.6
$ws.Range("H107").Value = 563.4286
$ws.Range("I107").Value = 509
$ws.Range("K107").Value = 1527
$ws.Range("M107").Value = 393
$ws.Range("H113").Value = 693.5454999999999
$ws.Range("I113").Value = 289.85715
$ws.Range("J113").Value = 1400
$ws.Range("K113").Value = 869.5714499999999
$ws.Range("L113").Value = 4200
$ws.Range("M113").Value = 1300.42855
$ws.Range("N113").Value = -8540
$ws.Range("H122").Value = 34484156
$ws.Range("I122").Value = 47620304
$ws.Range("J122").Value = 1768.25
$ws.Range("K122").Value = 142860912
$ws.Range("L122").Value = 5304.75
$ws.Range("M122").Value = -142858462
$ws.Range("N122").Value = -10204.75
$ws.Range("H126").Value = 7745.091
$ws.Range("I126").Value = 8419.6
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 25258.8
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -22788.8
$ws.Range("N126").Value = -7940
$ws.Range("H132").Value = 3658.45
$ws.Range("I132").Value = 3843.7144
$ws.Range("J132").Value = 2361.6
$ws.Range("K132").Value = 11531.1432
$ws.Range("L132").Value = 7084.799999999999
$ws.Range("M132").Value = -9001.143199999999
$ws.Range("N132").Value = -12144.8
